$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12 (item id 5515) on ALC
$ws.Range("H12").Value = 197.06667
$ws.Range("I12").Value = 204
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 204
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = -34
$ws.Range("N12").Value = -440

# Row 16 (item id 2146) on ALC
$ws.Range("H16").Value = 18653.637
$ws.Range("I16").Value = 6000
$ws.Range("J16").Value = 21465.555
$ws.Range("K16").Value = 6000
$ws.Range("L16").Value = 21465.555
$ws.Range("M16").Value = -5770
$ws.Range("N16").Value = -21925.555

# Row 18 (item id 5471) on ALC
$ws.Range("H18").Value = 225.28572
$ws.Range("I18").Value = 225.28572
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 225.28572
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 58.71428
$ws.Range("N18").ClearContents()

# Row 32 (item id 5484) on ALC
$ws.Range("H32").Value = 955.6667
$ws.Range("I32").Value = 475.25
$ws.Range("J32").Value = 1340
$ws.Range("K32").Value = 475.25
$ws.Range("L32").Value = 1340
$ws.Range("M32").Value = -149.25
$ws.Range("N32").Value = -1992

# Row 41 (item id 5478) on ALC
$ws.Range("H41").Value = 427.0909
$ws.Range("I41").Value = 466.33334
$ws.Range("J41").Value = 412.375
$ws.Range("K41").Value = 466.33334
$ws.Range("L41").Value = 412.375
$ws.Range("M41").Value = -26.33334000000002
$ws.Range("N41").Value = -1292.375

# Row 55 (item id 5517) on ALC
$ws.Range("H55").Value = 423.26666
$ws.Range("I55").Value = 87.8
$ws.Range("J55").Value = 591
$ws.Range("K55").Value = 87.8
$ws.Range("L55").Value = 591
$ws.Range("M55").Value = 126.2
$ws.Range("N55").Value = -1019

# Row 129 (item id 36115) on ALC
$ws.Range("H129").Value = 2536.087
$ws.Range("I129").Value = 478.5
$ws.Range("J129").Value = 4118.846
$ws.Range("K129").Value = 1435.5
$ws.Range("L129").Value = 12356.538
$ws.Range("M129").Value = 3564.5
$ws.Range("N129").Value = -22356.538

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (item id 27713) on ARM
$ws.Range("H2").Value = 822
$ws.Range("I2").Value = 822
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 822
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -709
$ws.Range("N2").ClearContents()

# Row 5 (item id 5091) on ARM
$ws.Range("H5").Value = 810.7857
$ws.Range("I5").Value = 1173.2222
$ws.Range("J5").Value = 158.4
$ws.Range("K5").Value = 1173.2222
$ws.Range("L5").Value = 158.4
$ws.Range("M5").Value = -1061.2222
$ws.Range("N5").Value = -382.4

# Row 32 (item id 44147) on ARM
$ws.Range("H32").Value = 27553.62
$ws.Range("I32").Value = 5438.825
$ws.Range("J32").Value = 116012.8
$ws.Range("K32").Value = 5438.825
$ws.Range("L32").Value = 116012.8
$ws.Range("M32").Value = -5151.825
$ws.Range("N32").Value = -116586.8

# Row 116 (item id 27713) on ARM
$ws.Range("H116").Value = 822
$ws.Range("I116").Value = 822
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 822
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1472
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (item id 27713) on BSM
$ws.Range("H3").Value = 822
$ws.Range("I3").Value = 822
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 822
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -708
$ws.Range("N3").ClearContents()

# Row 4 (item id 5091) on BSM
$ws.Range("H4").Value = 810.7857
$ws.Range("I4").Value = 1173.2222
$ws.Range("J4").Value = 158.4
$ws.Range("K4").Value = 1173.2222
$ws.Range("L4").Value = 158.4
$ws.Range("M4").Value = -1058.2222
$ws.Range("N4").Value = -388.4

# Row 5 (item id 1750) on BSM
$ws.Range("H5").Value = 1166.6666

# Row 22 (item id 5092) on BSM
$ws.Range("H22").Value = 2030.75
$ws.Range("I22").Value = 4173.8
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 4173.8
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -4000.8
$ws.Range("N22").Value = -846

# Row 128 (item id 38749) on BSM
$ws.Range("H128").Value = 946.6667
$ws.Range("I128").Value = 946.6667
$ws.Range("K128").Value = 2840.0001
$ws.Range("M128").Value = -350.0001000000002

# Row 134 (item id 43998) on BSM
$ws.Range("H134").Value = 8626461
$ws.Range("I134").Value = 10209104
$ws.Range("J134").Value = 9844.444
$ws.Range("K134").Value = 30627312
$ws.Range("L134").Value = 29533.332
$ws.Range("M134").Value = -30624777
$ws.Range("N134").Value = -34603.33199999999

$ws = $wb.Worksheets.Item("CUL")
# Row 33 (item id 4867) on CUL
$ws.Range("H33").Value = 2612.7856
$ws.Range("I33").Value = 999.5
$ws.Range("J33").Value = 3822.75
$ws.Range("K33").Value = 5997
$ws.Range("L33").Value = 22936.5
$ws.Range("M33").Value = -5714
$ws.Range("N33").Value = -23502.5

# Row 122 (item id 36078) on CUL
$ws.Range("H122").Value = 543.2439000000001
$ws.Range("I122").Value = 341.25
$ws.Range("J122").Value = 828.41174
$ws.Range("K122").Value = 3071.25
$ws.Range("L122").Value = 7455.70566
$ws.Range("M122").Value = -621.25
$ws.Range("N122").Value = -12355.70566

# Row 131 (item id 36060) on CUL
$ws.Range("H131").Value = 887.2033699999999
$ws.Range("I131").Value = 425.77777
$ws.Range("J131").Value = 970.26
$ws.Range("K131").Value = 1277.33331
$ws.Range("L131").Value = 2910.78
$ws.Range("M131").Value = 3762.66669
$ws.Range("N131").Value = -12990.78

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (item id 5277) on LTW
$ws.Range("H22").Value = 1359.4
$ws.Range("I22").Value = 1100.3334
$ws.Range("J22").Value = 1470.4286
$ws.Range("K22").Value = 1100.3334
$ws.Range("L22").Value = 1470.4286
$ws.Range("M22").Value = -805.3334
$ws.Range("N22").Value = -2060.4286

# Row 27 (item id 5277) on LTW
$ws.Range("H27").Value = 1359.4
$ws.Range("I27").Value = 1100.3334
$ws.Range("J27").Value = 1470.4286
$ws.Range("K27").Value = 1100.3334
$ws.Range("L27").Value = 1470.4286
$ws.Range("M27").Value = -993.3334
$ws.Range("N27").Value = -1684.4286

# Row 55 (item id 5284) on LTW
$ws.Range("H55").Value = 252.29033
$ws.Range("I55").Value = 275.94116
$ws.Range("J55").Value = 223.57143
$ws.Range("K55").Value = 275.94116
$ws.Range("L55").Value = 223.57143
$ws.Range("M55").Value = -102.94116
$ws.Range("N55").Value = -569.57143

# Row 132 (item id 44058) on LTW
$ws.Range("H132").Value = 316639.88
$ws.Range("I132").Value = 51966.9
$ws.Range("J132").Value = 669537.2
$ws.Range("K132").Value = 155900.7
$ws.Range("L132").Value = 2008611.6
$ws.Range("M132").Value = -153370.7
$ws.Range("N132").Value = -2013671.6
